$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value  = "stuendlich gemittelte CO-Konzentration"
$ws.Range("H3").Value  = "stuendlich gemittelte Sensorreaktion (nominell auf CO ausgerichtet) (Zinnoxid)"
$ws.Range("H4").Value  = "stuendlich gemittelte Gesamtkonzentration an nicht-metanischem Kohlenwasserstoff"
$ws.Range("H5").Value  = "stuendlich gemittelte Benzolkonzentration"
$ws.Range("H6").Value  = "stuendlich gemittelte Sensorreaktion (nominell auf NMHC ausgerichtet) (Titandioxid)"
$ws.Range("H7").Value  = "Echte stuendlich gemittelte NOx-Konzentration"
$ws.Range("H8").Value  = "stuendlich gemitteltes Sensoransprechverhalten (nominell auf NOx ausgerichtet)"
$ws.Range("H9").Value  = "stuendlich gemittelte NO2-Konzentration"
$ws.Range("H10").Value = "stuendlich gemittelte Sensorreaktion (nominell auf NO2 ausgerichtet) (Wolframoxid)"
$ws.Range("H11").Value = "stuendlich gemitteltes Sensoransprechverhalten (nominell O3-bezogen) (Indiumoxid)"
